$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: week number and covering-week date range
$ws.Range("A8").Value = "Volume 30   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# Fix up number/text formatting for cells whose type flips between
# a placeholder string ("0" / "***.*") and a real number, by copying
# the format from a same-column donor cell that already has the right style.
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("I16").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Copy() | Out-Null
$ws.Range("F26").PasteSpecial(-4122) | Out-Null
$ws.Range("I16").Copy() | Out-Null
$ws.Range("I26").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the new weekly crime-stat figures
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 1
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 106.25
$ws.Range("M16").Value = 10
$ws.Range("N16").Value = -81.967213114754
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "0"
$ws.Range("E17").Value = "***.*"
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -57.142857142857
$ws.Range("I17").Value = 27
$ws.Range("K17").Value = -42.553191489361
$ws.Range("L17").Value = 22.727272727272
$ws.Range("M17").Value = 58.823529411764
$ws.Range("N17").Value = -58.461538461538
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 14.285714285714
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -13.043478260869
$ws.Range("I18").Value = 64
$ws.Range("J18").Value = 83
$ws.Range("K18").Value = -22.891566265060
$ws.Range("L18").Value = -5.882352941176
$ws.Range("M18").Value = 25.490196078431
$ws.Range("N18").Value = -81.502890173410
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 24.390243902439
$ws.Range("I19").Value = 179
$ws.Range("J19").Value = 142
$ws.Range("K19").Value = 26.056338028169
$ws.Range("L19").Value = 67.289719626168
$ws.Range("M19").Value = 110.588235294118
$ws.Range("N19").Value = 70.476190476190
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 250
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 51.724137931034
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 37.5
$ws.Range("N20").Value = -81.967213114754
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = 2.061855670103
$ws.Range("I21").Value = 348
$ws.Range("J21").Value = 336
$ws.Range("K21").Value = 3.571428571428
$ws.Range("L21").Value = 46.835443037974
$ws.Range("M21").Value = 61.111111111111
$ws.Range("N21").Value = -63.174603174603
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 0
$ws.Range("C23").Value = "0"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = 11.111111111111
$ws.Range("L23").Value = 66.666666666666
$ws.Range("M23").Value = 11.111111111111
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 16
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = -17.333333333333
$ws.Range("I24").Value = 246
$ws.Range("J24").Value = 286
$ws.Range("K24").Value = -13.986013986014
$ws.Range("L24").Value = 24.242424242424
$ws.Range("M24").Value = 83.582089552238
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 61
$ws.Range("J25").Value = 75
$ws.Range("K25").Value = -18.666666666666
$ws.Range("L25").Value = 12.962962962963
$ws.Range("M25").Value = 12.962962962963
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 1
$ws.Range("K26").Value = -66.666666666666
$ws.Range("L26").Value = -80
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 142.857142857143
$ws.Range("L27").Value = 54.545454545454
